$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
# New column order: code, name, descr, bmtyp_code, lang_code, is_active,
#                    cr_by, cr_dtimes, upd_by, upd_dtimes, is_deleted, del_dtimes
$ws.Cells.Item(1,1).Value  = "code"
$ws.Cells.Item(1,2).Value  = "name"
$ws.Cells.Item(1,3).Value  = "descr"
$ws.Cells.Item(1,4).Value  = "bmtyp_code"
$ws.Cells.Item(1,5).Value  = "lang_code"
$ws.Cells.Item(1,6).Value  = "is_active"
$ws.Cells.Item(1,7).Value  = "cr_by"
$ws.Cells.Item(1,8).Value  = "cr_dtimes"
$ws.Cells.Item(1,9).Value  = "upd_by"
$ws.Cells.Item(1,10).Value = "upd_dtimes"
$ws.Cells.Item(1,11).Value = "is_deleted"
$ws.Cells.Item(1,12).Value = "del_dtimes"

# ---- Data rows (2-6) ----
$rows = @(
  @("TM","Les pouces","Empreinte des pouces gauche et droit","FNR"),
  @("RH","Dalle droite","Empreinte de la dalle droite","FNR"),
  @("LH","Dalle gauche","Empreinte de la dalle gauche","FNR"),
  @("LI","Iris gauche","Capture de l'iris Gauche","IRS"),
  @("RI","Iris droit","Capture de l'iris droit","IRS")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
  $r = $i + 2
  $row = $rows[$i]
  $ws.Cells.Item($r,1).Value = $row[0]
  $ws.Cells.Item($r,2).Value = $row[1]
  $ws.Cells.Item($r,3).Value = $row[2]
  $ws.Cells.Item($r,4).Value = $row[3]
  $ws.Cells.Item($r,5).Value = "fra"
  $ws.Cells.Item($r,6).Value = $true
  $ws.Cells.Item($r,7).Value = "superadmin"
  $ws.Cells.Item($r,8).Value = 45079.578128749999
  $ws.Cells.Item($r,8).NumberFormat = "mm:ss.0"
  $ws.Cells.Item($r,9).Value = "NULL"
  $ws.Cells.Item($r,10).Value = "NULL"
  $ws.Cells.Item($r,11).Value = $false
  $ws.Cells.Item($r,12).Value = "NULL"
}

# ---- Selection ----
$ws.Range("F16").Select() | Out-Null
